$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the existing hyperlinks; they'll be recreated below once the new
# row has been inserted and the URL cells are back in their final spots.
$ws.Hyperlinks.Delete()

# Insert a new row above the current row 2, pushing the two existing
# listings down to rows 3 and 4.
$ws.Rows.Item(2).Insert()

# Fill in the newly inserted row 2 with the new listing.
$ws.Range("A2").Value = "2025-11-01 12:33:54"
$ws.Range("B2").Value = "【募集】売上金額表と在庫管理のエクセル作成依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5425201"
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = "◇管理"

# The two pre-existing listings (now on rows 3 and 4) were re-scraped at
# the same time as the new one, so refresh their "retrieved at" stamp.
$ws.Range("A3").Value = "2025-11-01 12:33:54"
$ws.Range("A4").Value = "2025-11-01 12:33:54"

# Re-create the hyperlinks for all three URL cells, in row order, so the
# relationship ids line up the same way Excel would assign them.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5425201")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5425003")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5424906")

# Hyperlinks.Add re-styles the cell with a fresh (duplicate) style record;
# re-apply the named "Hyperlink" style so every URL cell shares the same
# style index again.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
